# Adds SPI Flash ROM (U3 / W25Q128) plus its supporting connector (J2) and
# two crystals (Q1, Q2) to the BOM worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Small in-place edits to existing rows
# ---------------------------------------------------------------------

# Row 3: crystal-load caps changed from 10pF (C1634) to 12pF (C38523)
$ws.Cells.Item(3, 4).Value = "12pF"
$ws.Cells.Item(3, 5).Value = "C38523"

# Row 5: decoupling-cap group grows by one (C18 added) -> qty 10 -> 11
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "C7, C8, C9, C10, C11, C12, C14, C15, C16, C17, C18"

# Row 10: LCSC part number update for the 4K7 resistors
$ws.Cells.Item(10, 5).Value = "C116693 "

# ---------------------------------------------------------------------
# 2) Row 16 becomes the new SPI flash ROM (U3) instead of the old J3 USB
#    connector (that connector is recreated below as J2 on row 17).
# ---------------------------------------------------------------------
$ws.Cells.Item(16, 2).Value = "U3"
$ws.Cells.Item(16, 3).Value = "SOIC-8_208mil"
$ws.Cells.Item(16, 4).Value = "W25Q128"
$ws.Cells.Item(16, 5).Value = "C97521"

# D16/E16 pick up word-wrap formatting (matches D15/E15 styling)
$ws.Range("D15:E16").WrapText = $true

# ---------------------------------------------------------------------
# 3) Row 17: USB micro-B connector, renamed J2 (was J3), qty 1
# ---------------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "J2"
$ws.Cells.Item(17, 3).Value = "10118192-0002LF"
$ws.Cells.Item(17, 4).Value = "USB-B-MICRO-SMD"
$ws.Cells.Item(17, 5).Value = "C2972784"
$ws.Cells.Item(17, 3).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$ws.Cells.Item(17, 4).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$ws.Cells.Item(17, 4).WrapText = $false
$ws.Cells.Item(17, 5).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral

# ---------------------------------------------------------------------
# 4) Row 18: Q1, 8 MHz crystal
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Q1"
$ws.Cells.Item(18, 3).Value = "HC-49S-SMD-2P"
$ws.Cells.Item(18, 4).Value = "8 Mhz Crystal 20pF"
$ws.Cells.Item(18, 5).Value = "C12674"
$ws.Cells.Item(18, 3).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$ws.Cells.Item(18, 4).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$ws.Cells.Item(18, 4).WrapText = $false
$ws.Cells.Item(18, 5).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$ws.Cells.Item(18, 5).Font.Bold = $false

# ---------------------------------------------------------------------
# 5) Row 19: Q2, 32.768 kHz crystal
# ---------------------------------------------------------------------
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Q2"
$ws.Cells.Item(19, 3).Value = "FC-12M"
$ws.Cells.Item(19, 4).Value = "32.768kHz Crystal 12.5pF"
$ws.Cells.Item(19, 5).Value = "C32346"
$ws.Cells.Item(19, 3).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$ws.Range("D19:E19").HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignGeneral
$ws.Range("D19:E19").WrapText = $true

# ---------------------------------------------------------------------
# 6) Row 20: SW1, SW2 tactile switches (previously row 17), unchanged
#    values but now row D picks up the "special" dark-grey font style.
# ---------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = 2
$ws.Cells.Item(20, 2).Value = "SW1, SW2"
$ws.Cells.Item(20, 3).Value = [char]0xFEFF + "EVQ-Q2"
$ws.Cells.Item(20, 4).Value = "SPST 6*6mm Tactile Switch"
$ws.Cells.Item(20, 5).Value = "C221880"
$ws.Cells.Item(20, 4).Font.Color = 1710618
$ws.Cells.Item(20, 4).HorizontalAlignment = [int][Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$ws.Cells.Item(20, 4).WrapText = $false

# ---------------------------------------------------------------------
# 7) Totals row moves from 19 (old formula "=SUM(A2:A16)") to 22 (row 21
#    stays blank) and now sums through the new rows. Row 19's old formula
#    content was already overwritten by the new Q2 data above, so all
#    that's left is to place the relocated total.
# ---------------------------------------------------------------------
$ws.Cells.Item(22, 1).Formula = "=SUM(A2:A20)"

$ws.Calculate()
